$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.788.84'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.635.83'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '215.43'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.257'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').Value = '19.70'
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '1.861.54'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').Value = '1.635.90'
$ws.Range('E14').Value = '  -1.98%  '
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '0.0₃0768'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '63.18'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').Value = '25.820.91'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').Value = '192.85'
$ws.Range('E21').Value = '  -1.37%  '
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').Value = '6.38'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = '1.82'
$ws.Range('E25').Value = '  +3.02%  '
$ws.Range('D26').Value = '142.16'
$ws.Range('E26').Value = '  +2.26%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('D37').Value = '1.132.55'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').Value = '0.542'
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('E42').Value = '  +0.70%  '
$ws.Range('D43').Value = '100.49'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '0.807'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').Value = '1.770.98'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('D47').Value = '55.33'
$ws.Range('E47').Value = '  -0.99%  '
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '2.32'
$ws.Range('E51').Value = '  +2.74%  '
